$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 682.94116
$ws.Range("I19").Value = 918.25
$ws.Range("J19").Value = 473.77777
$ws.Range("K19").Value = 918.25
$ws.Range("L19").Value = 473.77777
$ws.Range("M19").Value = -743.25
$ws.Range("N19").Value = -823.7777699999999
# Row 20
$ws.Range("H20").Value = 1300
$ws.Range("I20").Value = 1300
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1300
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1070
# Row 35
$ws.Range("H35").Value = 1300
$ws.Range("I35").Value = 1300
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1300
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -921
# Row 43
$ws.Range("H43").Value = 2784.5715
$ws.Range("I43").Value = 2914.3333
$ws.Range("J43").Value = 2687.25
$ws.Range("K43").Value = 2914.3333
$ws.Range("L43").Value = 2687.25
$ws.Range("M43").Value = -2845.3333
$ws.Range("N43").Value = -2825.25
# Row 80
$ws.Range("H80").Value = 609.8
$ws.Range("I80").Value = 649.125
$ws.Range("J80").Value = 583.5833
$ws.Range("K80").Value = 1947.375
$ws.Range("L80").Value = 1750.7499
$ws.Range("M80").Value = -949.375
$ws.Range("N80").Value = -3746.7499
# Row 83
$ws.Range("H83").Value = 609.8
$ws.Range("I83").Value = 649.125
$ws.Range("J83").Value = 583.5833
$ws.Range("K83").Value = 5842.125
$ws.Range("L83").Value = 5252.2497
$ws.Range("M83").Value = -850.125
$ws.Range("N83").Value = -15236.2497
# Row 116
$ws.Range("H116").Value = 3981.6155
$ws.Range("I116").Value = 3069.1428
$ws.Range("J116").Value = 5046.1665
$ws.Range("K116").Value = 3069.1428
$ws.Range("L116").Value = 5046.1665
$ws.Range("M116").Value = 372.8571999999999
$ws.Range("N116").Value = -11930.1665
# Row 132
$ws.Range("H132").Value = 2552.9285
$ws.Range("I132").Value = 2845.1667
$ws.Range("J132").Value = 799.5
$ws.Range("K132").Value = 8535.500100000001
$ws.Range("L132").Value = 2398.5
$ws.Range("M132").Value = -6005.500100000001
$ws.Range("N132").Value = -7458.5
# Row 138
$ws.Range("H138").Value = 2025.2
$ws.Range("I138").Value = 994.8889
$ws.Range("J138").Value = 2868.182
$ws.Range("K138").Value = 2984.6667
$ws.Range("L138").Value = 8604.545999999998
$ws.Range("M138").Value = 2155.3333
$ws.Range("N138").Value = -18884.546

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 2000
$ws.Range("I74").Value = 2000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1126
# Row 77
$ws.Range("H77").Value = 2000
$ws.Range("I77").Value = 2000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 10000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -5632
# Row 102
$ws.Range("H102").Value = 1830
$ws.Range("I102").Value = 1245
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1245
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 377
$ws.Range("N102").Value = -6244
# Row 132
$ws.Range("H132").Value = 600
$ws.Range("I132").Value = 600
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1800
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 730

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -77
# Row 86
$ws.Range("H86").Value = 2998.077
$ws.Range("I86").Value = 3139.2856
$ws.Range("J86").Value = 2833.3333
$ws.Range("K86").Value = 3139.2856
$ws.Range("L86").Value = 2833.3333
$ws.Range("M86").Value = -2016.2856
$ws.Range("N86").Value = -5079.3333
# Row 89
$ws.Range("H89").Value = 2998.077
$ws.Range("I89").Value = 3139.2856
$ws.Range("J89").Value = 2833.3333
$ws.Range("K89").Value = 15696.428
$ws.Range("L89").Value = 14166.6665
$ws.Range("M89").Value = -10080.428
$ws.Range("N89").Value = -25398.6665
# Row 107
$ws.Range("H107").Value = 1029.6666
$ws.Range("I107").Value = 886
$ws.Range("J107").Value = 1424.75
$ws.Range("K107").Value = 886
$ws.Range("L107").Value = 1424.75
$ws.Range("M107").Value = 1034
$ws.Range("N107").Value = -5264.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 9602.583000000001
$ws.Range("I86").Value = 10944.333
$ws.Range("J86").Value = 9155.333000000001
$ws.Range("K86").Value = 10944.333
$ws.Range("L86").Value = 9155.333000000001
$ws.Range("M86").Value = -9821.333000000001
$ws.Range("N86").Value = -11401.333
# Row 89
$ws.Range("H89").Value = 9602.583000000001
$ws.Range("I89").Value = 10944.333
$ws.Range("J89").Value = 9155.333000000001
$ws.Range("K89").Value = 54721.665
$ws.Range("L89").Value = 45776.665
$ws.Range("M89").Value = -49105.665
$ws.Range("N89").Value = -57008.665
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 1458.8
$ws.Range("I8").Value = 1458.8
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 4376.4
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -4237.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 83341170
$ws.Range("I70").Value = 166677840
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 166677840
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -166677570
$ws.Range("N70").Value = -5040
# Row 73
$ws.Range("H73").Value = 83341170
$ws.Range("I73").Value = 166677840
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 166677840
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -166676904
$ws.Range("N73").Value = -6372
# Row 102
$ws.Range("H102").Value = 1360.9
$ws.Range("I102").Value = 1345.4445
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1345.4445
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 276.5554999999999
$ws.Range("N102").Value = -4744
# Row 126
$ws.Range("H126").Value = 15256.143
$ws.Range("I126").Value = 10358.6
$ws.Range("J126").Value = 27500
$ws.Range("K126").Value = 31075.8
$ws.Range("L126").Value = 82500
$ws.Range("M126").Value = -28605.8
$ws.Range("N126").Value = -87440
# Row 132
$ws.Range("H132").Value = 1957.0834
$ws.Range("I132").Value = 1886.875
$ws.Range("J132").Value = 2097.5
$ws.Range("K132").Value = 5660.625
$ws.Range("L132").Value = 6292.5
$ws.Range("M132").Value = -3130.625
$ws.Range("N132").Value = -11352.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 16764.8
$ws.Range("I7").Value = 12449.714
$ws.Range("J7").Value = 26833.334
$ws.Range("K7").Value = 12449.714
$ws.Range("L7").Value = 26833.334
$ws.Range("M7").Value = -12337.714
$ws.Range("N7").Value = -27057.334
# Row 68
$ws.Range("H68").Value = 2466
$ws.Range("I68").Value = 2260.2
$ws.Range("J68").Value = 3495
$ws.Range("K68").Value = 2260.2
$ws.Range("L68").Value = 3495
$ws.Range("M68").Value = -1511.2
$ws.Range("N68").Value = -4993
# Row 71
$ws.Range("H71").Value = 2466
$ws.Range("I71").Value = 2260.2
$ws.Range("J71").Value = 3495
$ws.Range("K71").Value = 11301
$ws.Range("L71").Value = 17475
$ws.Range("M71").Value = -7557
$ws.Range("N71").Value = -24963
# Row 93
$ws.Range("H93").Value = 1643.875
$ws.Range("I93").Value = 1564.4286
$ws.Range("J93").Value = 2200
$ws.Range("K93").Value = 1564.4286
$ws.Range("L93").Value = 2200
$ws.Range("M93").Value = -316.4286
$ws.Range("N93").Value = -4696
# Row 126
$ws.Range("H126").Value = 16764.8
$ws.Range("I126").Value = 12449.714
$ws.Range("J126").Value = 26833.334
$ws.Range("K126").Value = 37349.142
$ws.Range("L126").Value = 80500.00199999999
$ws.Range("M126").Value = -34879.142
$ws.Range("N126").Value = -85440.00199999999
# Row 132
$ws.Range("H132").Value = 5083.5
$ws.Range("I132").Value = 5152.5
$ws.Range("J132").Value = 4738.5
$ws.Range("K132").Value = 15457.5
$ws.Range("L132").Value = 14215.5
$ws.Range("M132").Value = -12927.5
$ws.Range("N132").Value = -19275.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1066.6666
$ws.Range("I81").Value = 1050
$ws.Range("J81").Value = 1100
$ws.Range("K81").Value = 2100
$ws.Range("L81").Value = 2200
$ws.Range("M81").Value = -1039
$ws.Range("N81").Value = -4322
# Row 84
$ws.Range("H84").Value = 1066.6666
$ws.Range("I84").Value = 1050
$ws.Range("J84").Value = 1100
$ws.Range("K84").Value = 10500
$ws.Range("L84").Value = 11000
$ws.Range("M84").Value = -5196
$ws.Range("N84").Value = -21608
# Row 132
$ws.Range("H132").Value = 2534.4
$ws.Range("I132").Value = 2534.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7603.200000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5073.200000000001
# Row 136
$ws.Range("H136").Value = 5911.1113
$ws.Range("I136").Value = 6914.2856
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 20742.8568
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -18192.8568
$ws.Range("N136").Value = -12300
